$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a date-serial "last changed" value
# that is repeated down every data row (rows 2-227). Bump it from
# 46061 (2026-02-08) to 46062 (2026-02-09) for all rows.
$ws.Range("C2:C227").Value = 46062
